$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.898.63"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.21%  "
$ws.Range("D3").Value = "'1.765.23"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.73%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'329.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.47%  "
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").Value = "'0.4547"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.27%  "
$ws.Range("D8").Value = "'0.3515"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.45%  "
$ws.Range("D9").Value = "'41.97"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.39%  "
$ws.Range("D10").Value = "'0.07388"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.00%  "
$ws.Range("E11").Value = "  +1.31%  "
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("D13").Value = "'20.72"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.10%  "
$ws.Range("D14").Value = "'6.003"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.33%  "
$ws.Range("D15").Value = "'7.188"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.52%  "
$ws.Range("D16").Value = "'1.767.14"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.77%  "
$ws.Range("D17").Value = "'92.39"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.04%  "
$ws.Range("D18").Value = "'0.00001058"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.27%  "
$ws.Range("D19").Value = "'0.06442"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.90%  "
$ws.Range("E20").Value = "  +0.03%  "
$ws.Range("D21").Value = "'16.99"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.57%  "
$ws.Range("D22").Value = "'5.780"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.83%  "
$ws.Range("D23").Value = "'27.925.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.14%  "
$ws.Range("D24").Value = "'11.24"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.57%  "
$ws.Range("D25").Value = "'2.152"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.27%  "
$ws.Range("D26").Value = "'161.99"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.06%  "
$ws.Range("D27").Value = "'20.15"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("D28").Value = "'1.970.98"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.88%  "
$ws.Range("E29").Value = "  +2.96%  "
$ws.Range("D30").Value = "'123.90"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.28%  "
$ws.Range("D31").Value = "'1.076"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.77%  "
$ws.Range("D32").Value = "'0.09282"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.33%  "
$ws.Range("D33").Value = "'5.584"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.59%  "
$ws.Range("D34").Value = "'3.648"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.15%  "
$ws.Range("D35").Value = "'11.81"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.12%  "
$ws.Range("D36").Value = "'0.02273"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.35%  "
$ws.Range("D37").Value = "'0.06120"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.99%  "
$ws.Range("D38").Value = "'0.2087"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.15%  "
$ws.Range("D39").Value = "'4.940"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.51%  "
$ws.Range("D40").Value = "'0.6257"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.28%  "
$ws.Range("D41").Value = "'1.182"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.38%  "
$ws.Range("D42").Value = "'1.379"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.76%  "
$ws.Range("D43").Value = "'7.857"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.04%  "
$ws.Range("D44").Value = "'13.18"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.13%  "
$ws.Range("E45").Value = "  +0.45%  "
$ws.Range("D46").Value = "'0.5843"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.09%  "
$ws.Range("D47").Value = "'122.76"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'1.936"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.26%  "
$ws.Range("D49").Value = "'1.129"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("E50").Value = "  -1.11%  "
$ws.Range("D51").Value = "'72.91"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.20%  "
